$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, shifting existing rows 69-73 down to 70-74
$ws.Rows.Item(69).Insert()

# The newly inserted row 69 is blank; populate it with a full data record,
# matching the metadata columns shared by every row in this data block
# (same market/region/category info), then the new record's own values.
$ws.Cells.Item(69, 1).Value = $ws.Cells.Item(70, 1).Value2   # A Mercado ID
$ws.Cells.Item(69, 2).Value = $ws.Cells.Item(70, 2).Value2   # B Mercado
$ws.Cells.Item(69, 3).Value = $ws.Cells.Item(70, 3).Value2   # C Region
$ws.Cells.Item(69, 4).Value = 45124                          # D Fecha
$ws.Cells.Item(69, 4).NumberFormat = $ws.Cells.Item(70, 4).NumberFormat
$ws.Cells.Item(69, 5).Value = $ws.Cells.Item(70, 5).Value2   # E Codreg
$ws.Cells.Item(69, 6).Value = $ws.Cells.Item(70, 6).Value2   # F Categoria ID
$ws.Cells.Item(69, 7).Value = $ws.Cells.Item(70, 7).Value2   # G Categoria
$ws.Cells.Item(69, 8).Value = $ws.Cells.Item(70, 8).Value2   # H Variedad
$ws.Cells.Item(69, 9).Value = $ws.Cells.Item(70, 9).Value2   # I Calidad
$ws.Cells.Item(69, 10).Value = 10                            # J Volumen
$ws.Cells.Item(69, 11).Value = 12000                         # K Precio minimo
$ws.Cells.Item(69, 12).Value = 12000                         # L Precio maximo
$ws.Cells.Item(69, 13).Value = 12000                         # M Precio promedio ponderado
$ws.Cells.Item(69, 14).Value = $ws.Cells.Item(70, 14).Value2 # N Unidad de comercializacion
$ws.Cells.Item(69, 15).Value = $ws.Cells.Item(70, 15).Value2 # O Origen
$ws.Cells.Item(69, 16).Value = 1200                          # P Precio $/Kg
$ws.Cells.Item(69, 17).Value = $ws.Cells.Item(70, 17).Value2 # Q Kg o Unidades
$ws.Cells.Item(69, 18).Value = $ws.Cells.Item(70, 18).Value2 # R Clasificacion
